$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (MigrationInProgressAction) - reworded description text
$ws.Range("B6").Value = "What action should the robot take to the current file it is migrating. Options are either Move To In Progress (moves to a different folder before migrating) or No Action to keep in present location"

# Row 10 (MigrationSourceFileCleanUpCopyLocation) - reworded description text
$ws.Range("B10").Value = "The location to copy the source file to during file clean up operations. Only applies when MigrationSourceFileCleanUpAction is set to CopyToCompleted. Can use {0} in the path to insert the existing folder structure (ignoring the network drive it is on)."

# Row 11 (MigrationSourceFileCleanUpMoveLocation) - reworded description text
$ws.Range("B11").Value = "The location to move the source file to during file clean up operations. Only applies when MigrationSourceFileCleanUpAction is set to MoveToCompleted. Can use {0} in the path to insert the existing folder structure (ignoring the network drive it is on)."

# New row 15 - TargetFileExistsAction
$ws.Range("A15").Value = "TargetFileExistsAction"
$ws.Range("B15").Value = "If a file with the same name exists in the same location on the Target site, what action should be taken? Options are: Overwrite, Rename - Underscore Increment"
$ws.Range("C15").Value = "Rename - Underscore Increment"

# Update selection to match the committed state
$ws.Range("B7").Select()
